$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 250, shifting rows 250:294 down to 251:295
$ws.Rows.Item(250).Insert()

# Populate the new row 250 with the data from the diff
$ws.Cells.Item(250, 1).Value = 10
$ws.Cells.Item(250, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(250, 3).Value = "La Araucanía"
$ws.Cells.Item(250, 4).Value = 44637
$ws.Cells.Item(250, 5).Value = 9
$ws.Cells.Item(250, 6).Value = 100112009
$ws.Cells.Item(250, 7).Value = "Acelga"
$ws.Cells.Item(250, 8).Value = "Sin especificar"
$ws.Cells.Item(250, 9).Value = "Primera"
$ws.Cells.Item(250, 10).Value = 60
$ws.Cells.Item(250, 11).Value = 8000
$ws.Cells.Item(250, 12).Value = 8000
$ws.Cells.Item(250, 13).Value = 8000
$ws.Cells.Item(250, 14).Value = "`$/docena de atados (12 kilos)"
$ws.Cells.Item(250, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(250, 16).Value = 667
$ws.Cells.Item(250, 17).Value = 12
$ws.Cells.Item(250, 18).Value = "Hortaliza"
